$wb = $excel.ActiveWorkbook

# OFF sheet - update row 2 (Short Att, Short Comp, Deep Att, Deep Comp)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 278
$wsOff.Range("C2").Value = 195
$wsOff.Range("D2").Value = 62
$wsOff.Range("E2").Value = 22

# DEF sheet - update row 2 (Short Att, Short Comp, Deep Att, Deep Comp, Deep Int)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 285
$wsDef.Range("C2").Value = 201
$wsDef.Range("D2").Value = 55
$wsDef.Range("E2").Value = 18
$wsDef.Range("G2").Value = 4
